# "Alteração estilo Login e adicionado User"
# Append 5 new log rows (35-39) to the "ABR" sheet, continuing the existing
# numbering/pattern used by the table (A: incrementing counter, B: date text,
# C: fixed text, D: fixed quantity, E: machine code), matching the rows
# already present (r=2..34) but for date 10/03/2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABR")

$startRow = 35
$startCounter = 33

for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $counter = $startCounter + $i

    $ws.Cells.Item($row, 1).Value = $counter

    # Write "10/03/2022" as literal text (matching the other rows, which
    # store the date as plain text rather than a real date value). Assigning
    # the string directly would make Excel auto-convert it into a date
    # serial number, so we build it through a formula and convert the
    # formula's text result back into a plain value via copy / paste-values.
    $ws.Cells.Item($row, 2).Formula = '="10/03/2022"'
    $ws.Cells.Item($row, 2).Copy() | Out-Null
    $ws.Cells.Item($row, 2).PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($row, 3).Value = "ifjeifje"
    $ws.Cells.Item($row, 4).Value = 32
    $ws.Cells.Item($row, 5).Value = "cnc2"
}
